$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on target cells so values (including numeric-looking
# strings such as prices) are preserved exactly as text, matching the original
# inlineStr cell type rather than being auto-converted to numbers by Excel.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.212.69'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.300.05'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.37%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.81'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '103.36'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.79%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.610'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.61'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0911'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.35'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.45%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.967'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.36'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.651.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.311.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.373.34'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.40'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.68'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.60'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '280.75'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.07'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +16.43%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.21%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.88'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -4.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.39'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +6.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.95'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.06'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '164.27'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0876'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.87'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.136'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.62'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -10.16%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.55%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.59'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -2.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0349'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.76'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.78'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '100.50'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -4.22%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '69.34'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.227'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.04'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '112.00'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '77.12'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.31%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.30'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.605.68'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.72%  '
